$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.543.31'
$ws.Range("E2").Value = '  -1.35%  '
$ws.Range("D3").Value = '3.765.59'
$ws.Range("E3").Value = '  -2.05%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'595.62"
$ws.Range("E5").Value = '  -1.06%  '
$ws.Range("D6").Value = "'168.27"
$ws.Range("E6").Value = '  -0.67%  '
$ws.Range("D7").Value = '3.763.94'
$ws.Range("E7").Value = '  -2.13%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = "'0.523"
$ws.Range("E9").Value = '  -0.84%  '
$ws.Range("E10").Value = '  -2.05%  '
$ws.Range("D11").Value = "'6.45"
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("D12").Value = "'0.452"
$ws.Range("E12").Value = '  -1.50%  '
$ws.Range("D13").Value = "'0.0000273"
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("E14").Value = '  -2.30%  '
$ws.Range("D15").Value = '4.391.09'
$ws.Range("E15").Value = '  -2.21%  '
$ws.Range("D16").Value = '3.758.03'
$ws.Range("E16").Value = '  -1.86%  '
$ws.Range("D17").Value = "'18.54"
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = '67.445.36'
$ws.Range("E18").Value = '  -1.47%  '
$ws.Range("E19").Value = '  -3.29%  '
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("E21").Value = '  -5.20%  '
$ws.Range("D22").Value = "'466.58"
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("E23").Value = '  -2.64%  '
$ws.Range("D24").Value = "'83.51"
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("E25").Value = '  -9.46%  '
$ws.Range("E26").Value = '  -1.87%  '
$ws.Range("D27").Value = "'12.10"
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("D28").Value = "'10.25"
$ws.Range("E28").Value = '  +1.75%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  -2.38%  '
$ws.Range("D31").Value = '3.907.42'
$ws.Range("E31").Value = '  -2.10%  '
$ws.Range("E32").Value = '  -1.21%  '
$ws.Range("D33").Value = "'30.45"
$ws.Range("E33").Value = '  -3.80%  '
$ws.Range("E34").Value = '  -3.98%  '
$ws.Range("D35").Value = "'9.10"
$ws.Range("E35").Value = '  -3.05%  '
$ws.Range("D36").Value = '3.721.64'
$ws.Range("E36").Value = '  -2.22%  '
$ws.Range("D37").Value = "'3.80"
$ws.Range("E37").Value = '  +2.67%  '
$ws.Range("E38").Value = '  -1.32%  '
$ws.Range("E39").Value = '  -1.76%  '
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = '  -2.18%  '
$ws.Range("E41").Value = '  -3.11%  '
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").Value = '  -1.77%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = "'8.66"
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("D46").Value = "'1.93"
$ws.Range("E46").Value = '  -2.48%  '
$ws.Range("D47").Value = "'45.82"
$ws.Range("E47").Value = '  -2.83%  '
$ws.Range("D48").Value = "'394.88"
$ws.Range("E48").Value = '  -4.98%  '
$ws.Range("E49").Value = '  -7.22%  '
$ws.Range("D50").Value = "'139.20"
$ws.Range("E50").Value = '  -1.74%  '
$ws.Range("D51").Value = "'0.0352"
$ws.Range("E51").Value = '  -2.40%  '
